$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Swap the data (columns B..AC) of rows 200 and 201. Column A (the
#    sequential match id) stays where it is; only the match data moves.
# ---------------------------------------------------------------------
$colFirst = 2   # column B
$colLast  = 29  # column AC

for ($c = $colFirst; $c -le $colLast; $c++) {
    $cell200 = $ws.Cells.Item(200, $c)
    $cell201 = $ws.Cells.Item(201, $c)
    $v200 = $cell200.Value2
    $v201 = $cell201.Value2
    $cell200.Value2 = $v201
    $cell201.Value2 = $v200
}

# ---------------------------------------------------------------------
# 2) Update odds for the still-to-be-played match on row 223.
# ---------------------------------------------------------------------
$ws.Cells.Item(223, 14).Value2 = 1.166   # N223
$ws.Cells.Item(223, 15).Value2 = 7       # O223
$ws.Cells.Item(223, 16).Value2 = 15      # P223
$ws.Cells.Item(223, 17).Value2 = -2      # Q223
$ws.Cells.Item(223, 18).Value2 = 1.8     # R223
$ws.Cells.Item(223, 19).Value2 = 2       # S223
$ws.Cells.Item(223, 20).Value2 = 3.25    # T223

# ---------------------------------------------------------------------
# 3) Append two new fixture rows (224 and 225) with the same look and
#    feel (styles) as the existing data rows.
# ---------------------------------------------------------------------
$srcRange = $ws.Range("A223:AC223")
$dstRange = $ws.Range("A224:AC225")
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 224 ----------------------------------------------------------------
$ws.Cells.Item(224, 1).Value2  = 222          # A - id
$ws.Cells.Item(224, 2).Value2  = 8048492      # B
$ws.Cells.Item(224, 3).Value   = "Costa Rica Primera Division"   # C
$ws.Cells.Item(224, 4).Value   = "Costa Rica Primera Division"   # D
$ws.Cells.Item(224, 5).Value2  = 45387.95833333334                # E
$ws.Cells.Item(224, 6).Value   = "AD Guanacasteca"                # F
$ws.Cells.Item(224, 7).Value   = "Puntarenas"                     # G
$ws.Cells.Item(224, 11).Value2 = 2.2      # K
$ws.Cells.Item(224, 12).Value2 = 3.2      # L
$ws.Cells.Item(224, 13).Value2 = 3.3      # M
$ws.Cells.Item(224, 14).Value2 = 2.2      # N
$ws.Cells.Item(224, 15).Value2 = 3.2      # O
$ws.Cells.Item(224, 16).Value2 = 3.3      # P
$ws.Cells.Item(224, 17).Value2 = -0.25    # Q
$ws.Cells.Item(224, 18).Value2 = 1.85     # R
$ws.Cells.Item(224, 19).Value2 = 1.95     # S
$ws.Cells.Item(224, 20).Value2 = 2.25     # T
$ws.Cells.Item(224, 21).Value2 = 2        # U
$ws.Cells.Item(224, 22).Value2 = 1.8      # V
$ws.Cells.Item(224, 23).Value2 = 0        # W
$ws.Cells.Item(224, 24).Value2 = 0        # X
$ws.Cells.Item(224, 25).Value2 = 0        # Y
$ws.Cells.Item(224, 26).Value2 = 0        # Z
$ws.Cells.Item(224, 27).Value2 = 0        # AA

# clear the H/I/J/AB/AC cells (no format-pasted leftovers) for row 224
$ws.Cells.Item(224, 8).Clear()   | Out-Null
$ws.Cells.Item(224, 9).Clear()   | Out-Null
$ws.Cells.Item(224, 10).Clear()  | Out-Null
$ws.Cells.Item(224, 28).Clear()  | Out-Null
$ws.Cells.Item(224, 29).Clear()  | Out-Null

# Row 225 ----------------------------------------------------------------
$ws.Cells.Item(225, 1).Value2  = 223          # A - id
$ws.Cells.Item(225, 2).Value2  = 7623944      # B
$ws.Cells.Item(225, 3).Value   = "Costa Rica Primera Division"   # C
$ws.Cells.Item(225, 4).Value   = "Costa Rica Primera Division"   # D
$ws.Cells.Item(225, 5).Value2  = 45388.79166666666                # E
$ws.Cells.Item(225, 6).Value   = "Santos de Gupiles"               # F
$ws.Cells.Item(225, 7).Value   = "Municipal Liberia"                # G
$ws.Cells.Item(225, 11).Value2 = 2.9      # K
$ws.Cells.Item(225, 12).Value2 = 3.25     # L
$ws.Cells.Item(225, 13).Value2 = 2.375    # M
$ws.Cells.Item(225, 14).Value2 = 3.1      # N
$ws.Cells.Item(225, 15).Value2 = 3.3      # O
$ws.Cells.Item(225, 16).Value2 = 2.25     # P
$ws.Cells.Item(225, 17).Value2 = 0.25     # Q
$ws.Cells.Item(225, 18).Value2 = 1.85     # R
$ws.Cells.Item(225, 19).Value2 = 1.95     # S
$ws.Cells.Item(225, 20).Value2 = 2.5      # T
$ws.Cells.Item(225, 21).Value2 = 1.85     # U
$ws.Cells.Item(225, 22).Value2 = 1.95     # V
$ws.Cells.Item(225, 23).Value2 = 0        # W
$ws.Cells.Item(225, 24).Value2 = 0        # X
$ws.Cells.Item(225, 25).Value2 = 0        # Y
$ws.Cells.Item(225, 26).Value2 = 0        # Z
$ws.Cells.Item(225, 27).Value2 = 0        # AA

# clear the H/I/J/AB/AC cells (no format-pasted leftovers) for row 225
$ws.Cells.Item(225, 8).Clear()   | Out-Null
$ws.Cells.Item(225, 9).Clear()   | Out-Null
$ws.Cells.Item(225, 10).Clear()  | Out-Null
$ws.Cells.Item(225, 28).Clear()  | Out-Null
$ws.Cells.Item(225, 29).Clear()  | Out-Null

Write-Host "edit complete"
